# Fruta / hortaliza, semanal
# Prepend this week's two new price reports (rows 4 and 5) to the
# "Femacal de La Calera - Granada" data table, pushing all existing
# data rows down by two (old row 4 -> new row 6, ..., old row 34 -> new row 36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 4 (first data row after the header).
$ws.Rows("4:5").Insert()

# New row 4: Granada "Especial"
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Femacal de La Calera"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 45043
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100104
$ws.Range("H4").Value = "Frutos de pepita"
$ws.Range("I4").Value = 100104001
$ws.Range("J4").Value = "Granada"
$ws.Range("K4").Value = "Wonderfull"
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 45
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 17000
$ws.Range("Q4").Value = "`$/caja 14 kilos granel"
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 1214
$ws.Range("T4").Value = 14

# New row 5: Granada "Primera"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Femacal de La Calera"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 45043
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100104
$ws.Range("H5").Value = "Frutos de pepita"
$ws.Range("I5").Value = 100104001
$ws.Range("J5").Value = "Granada"
$ws.Range("K5").Value = "Wonderfull"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 67
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 14000
$ws.Range("Q5").Value = "`$/caja 14 kilos granel"
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 1000
$ws.Range("T5").Value = 14
